$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 129; this shifts the existing data in rows
# 129..235 down to rows 130..236 (extending the table by one record).
$ws.Rows.Item(129).Insert()

# Populate the newly inserted row 129 with the new price-report record.
$ws.Range("A129").Value = 5
$ws.Range("B129").Value = "Macroferia Regional de Talca"
$ws.Range("C129").Value = "Maule"
$ws.Range("D129").Value = 44669
$ws.Range("E129").Value = 7
$ws.Range("F129").Value = 100112024
$ws.Range("G129").Value = "Choclo"
$ws.Range("H129").Value = "Choclero"
$ws.Range("I129").Value = "Primera"
$ws.Range("J129").Value = 40000
$ws.Range("K129").Value = 180
$ws.Range("L129").Value = 180
$ws.Range("M129").Value = 180
$ws.Range("N129").Value = "`$/unidad"
$ws.Range("O129").Value = "Región del Maule"
$ws.Range("P129").Value = 180
$ws.Range("Q129").Value = 1
$ws.Range("R129").Value = "Hortaliza"
